$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy number formats from the (now shifted) column E into the new column D
# so the new column matches the formatting of its row neighbours.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the newest reporting period's figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 3100
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = 113800
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 162100
$ws.Range("D18").Value = -158900
$ws.Range("D20").Value = -5500
$ws.Range("D21").Value = -160900
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = -164400
$ws.Range("D24").Value = 600
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -165000
$ws.Range("D27").Value = -165000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 5500
$ws.Range("D33").Value = -165000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -165000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 456600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 100
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 9700
$ws.Range("D46").Value = 466400
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 18500
$ws.Range("D49").Value = 300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 3800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 489000
$ws.Range("D57").Value = 5100
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 22700
$ws.Range("D60").Value = 27700
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 69100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 96800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -291600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 392200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -165000
$ws.Range("D83").Value = 3500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -96200
$ws.Range("D91").Value = -2800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -2800
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 315900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 216900
$ws.Range("E24").Value = 1700
$ws.Range("E26").Value = -68400
$ws.Range("E27").Value = -68400
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
